$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = [char]39

$ws.Range("D2").Value = "42.898.72"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.213.81"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "$apos" + "256.21"
$ws.Range("E5").Value = "  +4.58%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "$apos" + "76.22"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "$apos" + "0.594"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").Value = "$apos" + "41.83"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").Value = "$apos" + "0.0907"
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").Value = "$apos" + "6.95"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "$apos" + "0.102"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "2.543.04"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "$apos" + "14.49"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "2.205.77"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").Value = "42.841.37"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").Value = "$apos" + "71.31"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "$apos" + "2.20"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "$apos" + "229.35"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "$apos" + "9.20"
$ws.Range("E24").Value = "  -8.37%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  +5.63%  "
$ws.Range("D27").Value = "$apos" + "10.73"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "$apos" + "3.33"
$ws.Range("E28").Value = "  -5.21%  "
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").Value = "$apos" + "2.19"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "$apos" + "174.15"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "$apos" + "0.0875"
$ws.Range("E32").Value = "  +9.68%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "$apos" + "20.29"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").Value = "$apos" + "0.0356"
$ws.Range("E36").Value = "  +7.40%  "
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("D38").Value = "$apos" + "4.29"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").Value = "$apos" + "12.60"
$ws.Range("E39").Value = "  -4.93%  "
$ws.Range("D40").Value = "$apos" + "2.82"
$ws.Range("E40").Value = "  +16.76%  "
$ws.Range("D41").Value = "$apos" + "2.11"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("E42").Value = "  -4.74%  "
$ws.Range("D43").Value = "$apos" + "0.198"
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("D44").Value = "$apos" + "60.05"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "$apos" + "102.40"
$ws.Range("E45").Value = "  -3.67%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "$apos" + "0.0976"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "$apos" + "8.29"
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("D48").Value = "$apos" + "0.456"
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "$apos" + "2.65"
$ws.Range("E51").Value = "  -0.87%  "
